$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Items")

# Add the new item rows (11-14) sketching the enemy AI "survival" data: Iron Armor (row 11),
# and three new stacked consumables - Herb, Mana Potion, Health Potion (rows 12-14).

# Row 11
$ws.Cells.Item(11, 1).Value = 0
$ws.Cells.Item(11, 2).Value = 3024
$ws.Cells.Item(11, 3).Value = 1010
$ws.Cells.Item(11, 4).Value = "a"
$ws.Cells.Item(11, 5).Value = "p"
$ws.Cells.Item(11, 6).Value = "h"
$ws.Cells.Item(11, 7).Value = "c"
$ws.Cells.Item(11, 8).Value = 50
$ws.Cells.Item(11, 9).Value = 255
$ws.Cells.Item(11, 10).Value = 0
$ws.Cells.Item(11, 11).Value = 255
$ws.Cells.Item(11, 12).Value = "Iron Armor"
$ws.Cells.Item(11, 13).Value = 1
$ws.Cells.Item(11, 14).Value = 3
$ws.Cells.Item(11, 15).Value = 0
$ws.Cells.Item(11, 16).Value = 0
$ws.Cells.Item(11, 17).Value = 0
$ws.Cells.Item(11, 18).Value = 0
$ws.Cells.Item(11, 19).Value = 0
$ws.Cells.Item(11, 20).Value = 0
$ws.Cells.Item(11, 21).Value = 0
$ws.Cells.Item(11, 22).Value = 0
$ws.Cells.Item(11, 23).Value = 0
$ws.Cells.Item(11, 24).Value = 0
$ws.Cells.Item(11, 25).Value = 0
$ws.Cells.Item(11, 26).Value = 0
$ws.Cells.Item(11, 27).Value = 0
$ws.Cells.Item(11, 28).Value = 0
$ws.Cells.Item(11, 29).Value = 0
$ws.Cells.Item(11, 30).Value = 0
$ws.Cells.Item(11, 31).Value = 0
$ws.Cells.Item(11, 32).Value = 0
$ws.Cells.Item(11, 33).Value = 0
$ws.Cells.Item(11, 34).Value = 0
$ws.Cells.Item(11, 35).Value = 0
$ws.Cells.Item(11, 36).Value = 0
$ws.Cells.Item(11, 37).Value = 2
$ws.Cells.Item(11, 38).Value = 2
$ws.Cells.Item(11, 39).Value = 2
$ws.Cells.Item(11, 40).Value = 2
$ws.Cells.Item(11, 41).Value = 0
$ws.Cells.Item(11, 42).Value = 0
$ws.Cells.Item(11, 43).Value = 0
$ws.Cells.Item(11, 44).Value = 0
$ws.Cells.Item(11, 45).Value = 0
$ws.Cells.Item(11, 46).Value = 0
$ws.Cells.Item(11, 47).Value = 0
$ws.Cells.Item(11, 48).Value = 0
$ws.Cells.Item(11, 49).Value = 0
$ws.Cells.Item(11, 50).Value = 1
$ws.Cells.Item(11, 51).Value = 0
$ws.Cells.Item(11, 52).Value = "1,10,-1"
$ws.Cells.Item(11, 53).Value = "2,30,30,-1"
$ws.Cells.Item(11, 54).Value = "7,5,5,5,5,5,5,5,3,7"
$ws.Cells.Item(11, 55).Value = -1
$ws.Cells.Item(11, 56).Value = -1
$ws.Cells.Item(11, 57).Value = -1
$ws.Cells.Item(11, 58).Value = -1
$ws.Cells.Item(11, 59).Value = -1
$ws.Cells.Item(11, 60).Value = "13,10,10,10,10,10,10,10,10,10,10,10,10,100,-1"
$ws.Cells.Item(11, 61).Value = -1
$ws.Cells.Item(11, 62).Value = -1
$ws.Cells.Item(11, 63).Value = "A sturdy suit of heavy armor.&&(+2 BluntDR, +2 ChopDr, +2 PierceDR, +2 SlashDR)"

# Row 12
$ws.Cells.Item(12, 1).Value = 200
$ws.Cells.Item(12, 2).Value = 3021
$ws.Cells.Item(12, 3).Value = 1011
$ws.Cells.Item(12, 4).Value = "i"
$ws.Cells.Item(12, 5).Value = "p"
$ws.Cells.Item(12, 6).Value = "l"
$ws.Cells.Item(12, 7).Value = "d"
$ws.Cells.Item(12, 8).Value = 10
$ws.Cells.Item(12, 9).Value = 255
$ws.Cells.Item(12, 10).Value = 0
$ws.Cells.Item(12, 11).Value = 255
$ws.Cells.Item(12, 12).Value = "Herb"
$ws.Cells.Item(12, 13).Value = 3
$ws.Cells.Item(12, 14).Value = 3
$ws.Cells.Item(12, 15).Value = 0
$ws.Cells.Item(12, 16).Value = 0
$ws.Cells.Item(12, 17).Value = 0
$ws.Cells.Item(12, 18).Value = 0
$ws.Cells.Item(12, 19).Value = 0
$ws.Cells.Item(12, 20).Value = 0
$ws.Cells.Item(12, 21).Value = 0
$ws.Cells.Item(12, 22).Value = 0
$ws.Cells.Item(12, 23).Value = 0
$ws.Cells.Item(12, 24).Value = 0
$ws.Cells.Item(12, 25).Value = 2
$ws.Cells.Item(12, 26).Value = 0
$ws.Cells.Item(12, 27).Value = 0
$ws.Cells.Item(12, 28).Value = 0
$ws.Cells.Item(12, 29).Value = 0
$ws.Cells.Item(12, 30).Value = 0
$ws.Cells.Item(12, 31).Value = 0
$ws.Cells.Item(12, 32).Value = 1
$ws.Cells.Item(12, 33).Value = 2
$ws.Cells.Item(12, 34).Value = 5
$ws.Cells.Item(12, 35).Value = 0
$ws.Cells.Item(12, 36).Value = 0
$ws.Cells.Item(12, 37).Value = 0
$ws.Cells.Item(12, 38).Value = 0
$ws.Cells.Item(12, 39).Value = 0
$ws.Cells.Item(12, 40).Value = 0
$ws.Cells.Item(12, 41).Value = 0
$ws.Cells.Item(12, 42).Value = 0
$ws.Cells.Item(12, 43).Value = 0
$ws.Cells.Item(12, 44).Value = 0
$ws.Cells.Item(12, 45).Value = 0
$ws.Cells.Item(12, 46).Value = 0
$ws.Cells.Item(12, 47).Value = 0
$ws.Cells.Item(12, 48).Value = 0
$ws.Cells.Item(12, 49).Value = 0
$ws.Cells.Item(12, 50).Value = 1
$ws.Cells.Item(12, 51).Value = 0
$ws.Cells.Item(12, 52).Value = "1,10,-1"
$ws.Cells.Item(12, 53).Value = -1
$ws.Cells.Item(12, 54).Value = -1
$ws.Cells.Item(12, 55).Value = -1
$ws.Cells.Item(12, 56).Value = -1
$ws.Cells.Item(12, 57).Value = -1
$ws.Cells.Item(12, 58).Value = -1
$ws.Cells.Item(12, 59).Value = -1
$ws.Cells.Item(12, 60).Value = -1
$ws.Cells.Item(12, 60).NumberFormat = "#,##0"
$ws.Cells.Item(12, 61).Value = -1
$ws.Cells.Item(12, 62).Value = -1
$ws.Cells.Item(12, 63).Value = "A bag of health restoring herbs.&&(+1 Dam, 2-5 Turns, Restores 2 HP per turn)"

# Row 13
$ws.Cells.Item(13, 1).Value = 200
$ws.Cells.Item(13, 2).Value = 3023
$ws.Cells.Item(13, 3).Value = 1012
$ws.Cells.Item(13, 4).Value = "i"
$ws.Cells.Item(13, 5).Value = "p"
$ws.Cells.Item(13, 6).Value = "l"
$ws.Cells.Item(13, 7).Value = "c"
$ws.Cells.Item(13, 8).Value = 20
$ws.Cells.Item(13, 9).Value = 255
$ws.Cells.Item(13, 10).Value = 0
$ws.Cells.Item(13, 11).Value = 255
$ws.Cells.Item(13, 12).Value = "Mana Potion"
$ws.Cells.Item(13, 13).Value = 1
$ws.Cells.Item(13, 14).Value = 2
$ws.Cells.Item(13, 15).Value = 0
$ws.Cells.Item(13, 16).Value = 0
$ws.Cells.Item(13, 17).Value = 0
$ws.Cells.Item(13, 18).Value = 0
$ws.Cells.Item(13, 19).Value = 0
$ws.Cells.Item(13, 20).Value = 0
$ws.Cells.Item(13, 21).Value = 0
$ws.Cells.Item(13, 22).Value = 0
$ws.Cells.Item(13, 23).Value = 0
$ws.Cells.Item(13, 24).Value = 0
$ws.Cells.Item(13, 25).Value = 0
$ws.Cells.Item(13, 26).Value = 0
$ws.Cells.Item(13, 27).Value = 6
$ws.Cells.Item(13, 28).Value = 0
$ws.Cells.Item(13, 29).Value = 0
$ws.Cells.Item(13, 30).Value = 0
$ws.Cells.Item(13, 31).Value = 0
$ws.Cells.Item(13, 32).Value = 0
$ws.Cells.Item(13, 33).Value = 0
$ws.Cells.Item(13, 34).Value = 0
$ws.Cells.Item(13, 35).Value = 0
$ws.Cells.Item(13, 36).Value = 0
$ws.Cells.Item(13, 37).Value = 0
$ws.Cells.Item(13, 38).Value = 0
$ws.Cells.Item(13, 39).Value = 0
$ws.Cells.Item(13, 40).Value = 0
$ws.Cells.Item(13, 41).Value = 0
$ws.Cells.Item(13, 42).Value = 0
$ws.Cells.Item(13, 43).Value = 0
$ws.Cells.Item(13, 44).Value = 0
$ws.Cells.Item(13, 45).Value = 0
$ws.Cells.Item(13, 46).Value = 0
$ws.Cells.Item(13, 47).Value = 0
$ws.Cells.Item(13, 48).Value = 0
$ws.Cells.Item(13, 49).Value = 0
$ws.Cells.Item(13, 50).Value = 1
$ws.Cells.Item(13, 51).Value = 0
$ws.Cells.Item(13, 52).Value = "1,10,-1"
$ws.Cells.Item(13, 53).Value = -1
$ws.Cells.Item(13, 54).Value = -1
$ws.Cells.Item(13, 55).Value = -1
$ws.Cells.Item(13, 56).Value = -1
$ws.Cells.Item(13, 57).Value = -1
$ws.Cells.Item(13, 58).Value = -1
$ws.Cells.Item(13, 59).Value = -1
$ws.Cells.Item(13, 60).Value = -1
$ws.Cells.Item(13, 60).NumberFormat = "#,##0"
$ws.Cells.Item(13, 61).Value = -1
$ws.Cells.Item(13, 62).Value = -1
$ws.Cells.Item(13, 63).Value = "A potion which restores mana.&&(Restores 4 Mana)"

# Row 14
$ws.Cells.Item(14, 1).Value = 200
$ws.Cells.Item(14, 2).Value = 3022
$ws.Cells.Item(14, 3).Value = 1013
$ws.Cells.Item(14, 4).Value = "i"
$ws.Cells.Item(14, 5).Value = "p"
$ws.Cells.Item(14, 6).Value = "l"
$ws.Cells.Item(14, 7).Value = "c"
$ws.Cells.Item(14, 8).Value = 15
$ws.Cells.Item(14, 9).Value = 255
$ws.Cells.Item(14, 10).Value = 0
$ws.Cells.Item(14, 11).Value = 255
$ws.Cells.Item(14, 12).Value = "Health Potion"
$ws.Cells.Item(14, 13).Value = 2
$ws.Cells.Item(14, 14).Value = 3
$ws.Cells.Item(14, 15).Value = 0
$ws.Cells.Item(14, 16).Value = 0
$ws.Cells.Item(14, 17).Value = 0
$ws.Cells.Item(14, 18).Value = 0
$ws.Cells.Item(14, 19).Value = 0
$ws.Cells.Item(14, 20).Value = 0
$ws.Cells.Item(14, 21).Value = 0
$ws.Cells.Item(14, 22).Value = 0
$ws.Cells.Item(14, 23).Value = 0
$ws.Cells.Item(14, 24).Value = 0
$ws.Cells.Item(14, 25).Value = 10
$ws.Cells.Item(14, 26).Value = 0
$ws.Cells.Item(14, 27).Value = 0
$ws.Cells.Item(14, 28).Value = 0
$ws.Cells.Item(14, 29).Value = 0
$ws.Cells.Item(14, 30).Value = 0
$ws.Cells.Item(14, 31).Value = 0
$ws.Cells.Item(14, 32).Value = 0
$ws.Cells.Item(14, 33).Value = 0
$ws.Cells.Item(14, 34).Value = 0
$ws.Cells.Item(14, 35).Value = 0
$ws.Cells.Item(14, 36).Value = 0
$ws.Cells.Item(14, 37).Value = 0
$ws.Cells.Item(14, 38).Value = 0
$ws.Cells.Item(14, 39).Value = 0
$ws.Cells.Item(14, 40).Value = 0
$ws.Cells.Item(14, 41).Value = 0
$ws.Cells.Item(14, 42).Value = 0
$ws.Cells.Item(14, 43).Value = 0
$ws.Cells.Item(14, 44).Value = 0
$ws.Cells.Item(14, 45).Value = 0
$ws.Cells.Item(14, 46).Value = 0
$ws.Cells.Item(14, 47).Value = 0
$ws.Cells.Item(14, 48).Value = 0
$ws.Cells.Item(14, 49).Value = 0
$ws.Cells.Item(14, 50).Value = 1
$ws.Cells.Item(14, 51).Value = 0
$ws.Cells.Item(14, 52).Value = "1,10,-1"
$ws.Cells.Item(14, 53).Value = -1
$ws.Cells.Item(14, 54).Value = -1
$ws.Cells.Item(14, 55).Value = -1
$ws.Cells.Item(14, 56).Value = -1
$ws.Cells.Item(14, 57).Value = -1
$ws.Cells.Item(14, 58).Value = -1
$ws.Cells.Item(14, 59).Value = -1
$ws.Cells.Item(14, 60).Value = -1
$ws.Cells.Item(14, 60).NumberFormat = "#,##0"
$ws.Cells.Item(14, 61).Value = -1
$ws.Cells.Item(14, 62).Value = -1
$ws.Cells.Item(14, 63).Value = "A potion which restores health.&&(Restores 10 HP)"

# Update the sheet view: scroll back to the left edge and select F12 (where we were last editing)
$ws.Activate()
$ws.Range("F12").Select()
